$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("N96").ClearContents()
$ws.Range("H96").Value = 513
$ws.Range("I96").Value = 513
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 1539
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -166
$ws.Range("H132").Value = 4275455.5
$ws.Range("I132").Value = 1565.3334
$ws.Range("J132").Value = 37041948
$ws.Range("K132").Value = 4696.0002
$ws.Range("L132").Value = 111125844
$ws.Range("M132").Value = -2166.0002
$ws.Range("N132").Value = -111130904
$ws.Range("H137").Value = 1386.8928
$ws.Range("I137").Value = 1432.6842
$ws.Range("J137").Value = 1290.2222
$ws.Range("K137").Value = 4298.0526
$ws.Range("L137").Value = 3870.6666
$ws.Range("M137").Value = -1748.0526
$ws.Range("N137").Value = -8970.6666
$ws.Range("H138").Value = 3299.831
$ws.Range("I138").Value = 1694.8
$ws.Range("J138").Value = 4324.3193
$ws.Range("K138").Value = 5084.4
$ws.Range("L138").Value = 12972.9579
$ws.Range("M138").Value = 55.60000000000036
$ws.Range("N138").Value = -23252.9579

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 201000.2
$ws.Range("I2").Value = 223000.22
$ws.Range("K2").Value = 223000.22
$ws.Range("M2").Value = -222887.22
$ws.Range("H32").Value = 16230.444
$ws.Range("I32").Value = 11686.844
$ws.Range("J32").Value = 27414.691
$ws.Range("K32").Value = 11686.844
$ws.Range("L32").Value = 27414.691
$ws.Range("M32").Value = -11399.844
$ws.Range("N32").Value = -27988.691
$ws.Range("H61").Value = 503304.4
$ws.Range("I61").Value = 3125.2666
$ws.Range("K61").Value = 3125.2666
$ws.Range("M61").Value = -2913.2666
$ws.Range("H116").Value = 201000.2
$ws.Range("I116").Value = 223000.22
$ws.Range("K116").Value = 223000.22
$ws.Range("M116").Value = -220706.22
$ws.Range("H134").Value = 65907.5
$ws.Range("J134").Value = 65907.5
$ws.Range("L134").Value = 65907.5
$ws.Range("N134").Value = -76047.5
$ws.Range("H136").Value = 503304.4
$ws.Range("I136").Value = 3125.2666
$ws.Range("K136").Value = 9375.799800000001
$ws.Range("M136").Value = -6825.799800000001

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 201000.2
$ws.Range("I3").Value = 223000.22
$ws.Range("K3").Value = 223000.22
$ws.Range("M3").Value = -222886.22

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("M31").ClearContents()
$ws.Range("H31").Value = 15159524
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 15159524
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 15159524
$ws.Range("N31").Value = -15160114
$ws.Range("M34").ClearContents()
$ws.Range("H34").Value = 15159524
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 15159524
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 15159524
$ws.Range("N34").Value = -15159928
$ws.Range("M38").ClearContents()
$ws.Range("N38").ClearContents()
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("N46").ClearContents()
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 0

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 5005118.5
$ws.Range("J33").Value = 9202.727999999999
$ws.Range("L33").Value = 55216.36799999999
$ws.Range("N33").Value = -55782.36799999999
$ws.Range("N99").ClearContents()
$ws.Range("H99").Value = 2925
$ws.Range("I99").Value = 2925
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 8775
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -6529
$ws.Range("H108").Value = 50463.5
$ws.Range("I108").Value = 50463.5
$ws.Range("K108").Value = 151390.5
$ws.Range("M108").Value = -148510.5
$ws.Range("H109").Value = 1698.5714
$ws.Range("I109").Value = 868.5714
$ws.Range("J109").Value = 2528.5715
$ws.Range("K109").Value = 2605.7142
$ws.Range("L109").Value = 7585.7145
$ws.Range("M109").Value = -1565.7142
$ws.Range("N109").Value = -9665.7145
$ws.Range("H125").Value = 5458.615
$ws.Range("I125").Value = 3746.2
$ws.Range("J125").Value = 6528.875
$ws.Range("K125").Value = 11238.6
$ws.Range("L125").Value = 19586.625
$ws.Range("M125").Value = -6318.599999999999
$ws.Range("N125").Value = -29426.625
$ws.Range("H131").Value = 2128706
$ws.Range("I131").Value = 6667356.5
$ws.Range("J131").Value = 1213.3125
$ws.Range("K131").Value = 20002069.5
$ws.Range("L131").Value = 3639.9375
$ws.Range("M131").Value = -19997029.5
$ws.Range("N131").Value = -13719.9375
$ws.Range("H134").Value = 11930
$ws.Range("I134").Value = 14663.75
$ws.Range("J134").Value = 9500
$ws.Range("K134").Value = 43991.25
$ws.Range("L134").Value = 28500
$ws.Range("M134").Value = -38921.25
$ws.Range("N134").Value = -38640
$ws.Range("H139").Value = 5285.793
$ws.Range("I139").Value = 6682.778
$ws.Range("J139").Value = 2999.818
$ws.Range("K139").Value = 20048.334
$ws.Range("L139").Value = 8999.454000000002
$ws.Range("M139").Value = -14908.334
$ws.Range("N139").Value = -19279.454
$ws.Range("H140").Value = 2378.8333
$ws.Range("I140").Value = 2378.8333
$ws.Range("K140").Value = 7136.499899999999
$ws.Range("M140").Value = -1956.499899999999

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 13000000
$ws.Range("I7").Value = 13000000
$ws.Range("K7").Value = 13000000
$ws.Range("M7").Value = -12999888
$ws.Range("H8").Value = 13000000
$ws.Range("I8").Value = 13000000
$ws.Range("K8").Value = 13000000
$ws.Range("M8").Value = -12999861

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 15855.615
$ws.Range("I136").Value = 7400.8
$ws.Range("J136").Value = 21139.875
$ws.Range("K136").Value = 22202.4
$ws.Range("L136").Value = 63419.625
$ws.Range("M136").Value = -19652.4
$ws.Range("N136").Value = -68519.625

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 12827129
$ws.Range("I136").Value = 11241.6
$ws.Range("J136").Value = 20837060
$ws.Range("K136").Value = 33724.8
$ws.Range("L136").Value = 62511180
$ws.Range("M136").Value = -31174.8
$ws.Range("N136").Value = -62516280
